# "Fixed and updated the OneStopShop"
#
# Summary of the edit being replayed:
#   1. Duplicate the "aggregate_data" sheet, place the copy at the end of the
#      workbook and rename it "aggregate_data1" (keeps the original data /
#      old formatting, exactly like the sheet it was copied from before the
#      subsequent formatting tweak below).
#   2. Normalize the stale/narrow sheet formatting on "taken1" and
#      "aggregate_data" (old 8-pt/15-pt defaults) to match the rest of the
#      workbook's current defaults (10-pt base col width / 16-pt row
#      height).
#   3. Re-point the workbook's active tab/selection at Sheet1 (instead of
#      "taken"), scrolled so column F is visible, with B6 selected.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "aggregate_data" -> "aggregate_data1" ----------------
$src  = $wb.Worksheets.Item("aggregate_data")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "aggregate_data1"

# --- 2. Bring "taken1" / "aggregate_data" formatting up to date --------
$taken1 = $wb.Worksheets.Item("taken1")
$taken1.StandardWidth = 8.83203125
$taken1.StandardHeight = 16

$aggregateData = $wb.Worksheets.Item("aggregate_data")
$aggregateData.StandardWidth = 8.83203125
$aggregateData.StandardHeight = 16

# --- 3. Move the active tab/selection back to Sheet1 --------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Activate()
[void]$sheet1.Range("B6").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
